# Tutorial 6 solution update: switch the "Date" column from DD/MM/YYYY to
# DD-MM-YYYY formatting for every attendance row, and flip the
# Total-Attendance-Count (D3) / Invalid (G3) counters for the first date
# row from 0 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new date text (slashes replaced with dashes).
$dates = [ordered]@{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

foreach ($row in $dates.Keys) {
    $text = $dates[$row]
    $cell = $ws.Cells.Item($row, 1)

    # Dates whose day-of-month is <= 12 read as ambiguous (day-first vs
    # month-first), so a plain assignment would get auto-recognised and
    # converted into a date serial number instead of staying literal text.
    # A leading apostrophe forces text entry; resetting the style back to
    # Normal afterwards drops the "quote prefix" formatting so the cell
    # ends up identical to a plain text write (still just a string, no
    # extra number formatting applied).
    $day = [int]($text.Split("-")[0])
    if ($day -le 12) {
        $cell.Value = "'" + $text
        $cell.Style = "Normal"
    } else {
        $cell.Value = $text
    }
}

# Row 3 (28-07-2022): Total Attendance Count and Invalid both go 0 -> 1.
$ws.Cells.Item(3, 4).Value = 1
$ws.Cells.Item(3, 7).Value = 1
